$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 fresh weekly rows above the existing data block (old row 891
# onward shifts down to 895+), growing the used range from R953 to R957.
$ws.Rows("891:894").Insert()

# Row 891: Coliflor, Primera, Región Metropolitana
$ws.Cells.Item(891, 1).Value = 9
$ws.Cells.Item(891, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(891, 3).Value = "Metropolitana"
$ws.Cells.Item(891, 4).Value = 45106
$ws.Cells.Item(891, 5).Value = 13
$ws.Cells.Item(891, 6).Value = 100112008
$ws.Cells.Item(891, 7).Value = "Coliflor"
$ws.Cells.Item(891, 8).Value = "Sin especificar"
$ws.Cells.Item(891, 9).Value = "Primera"
$ws.Cells.Item(891, 10).Value = 1600
$ws.Cells.Item(891, 11).Value = 800
$ws.Cells.Item(891, 12).Value = 900
$ws.Cells.Item(891, 13).Value = 850
$ws.Cells.Item(891, 14).Value = "$/unidad"
$ws.Cells.Item(891, 15).Value = "Región Metropolitana"
$ws.Cells.Item(891, 16).Value = 850
$ws.Cells.Item(891, 17).Value = 1
$ws.Cells.Item(891, 18).Value = "Hortaliza"

# Row 892: Coliflor, Primera, Región de O'Higgins
$ws.Cells.Item(892, 1).Value = 9
$ws.Cells.Item(892, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(892, 3).Value = "Metropolitana"
$ws.Cells.Item(892, 4).Value = 45106
$ws.Cells.Item(892, 5).Value = 13
$ws.Cells.Item(892, 6).Value = 100112008
$ws.Cells.Item(892, 7).Value = "Coliflor"
$ws.Cells.Item(892, 8).Value = "Sin especificar"
$ws.Cells.Item(892, 9).Value = "Primera"
$ws.Cells.Item(892, 10).Value = 970
$ws.Cells.Item(892, 11).Value = 900
$ws.Cells.Item(892, 12).Value = 1000
$ws.Cells.Item(892, 13).Value = 950
$ws.Cells.Item(892, 14).Value = "$/unidad"
$ws.Cells.Item(892, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(892, 16).Value = 950
$ws.Cells.Item(892, 17).Value = 1
$ws.Cells.Item(892, 18).Value = "Hortaliza"

# Row 893: Coliflor, Segunda, Región Metropolitana
$ws.Cells.Item(893, 1).Value = 9
$ws.Cells.Item(893, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(893, 3).Value = "Metropolitana"
$ws.Cells.Item(893, 4).Value = 45106
$ws.Cells.Item(893, 5).Value = 13
$ws.Cells.Item(893, 6).Value = 100112008
$ws.Cells.Item(893, 7).Value = "Coliflor"
$ws.Cells.Item(893, 8).Value = "Sin especificar"
$ws.Cells.Item(893, 9).Value = "Segunda"
$ws.Cells.Item(893, 10).Value = 970
$ws.Cells.Item(893, 11).Value = 700
$ws.Cells.Item(893, 12).Value = 700
$ws.Cells.Item(893, 13).Value = 700
$ws.Cells.Item(893, 14).Value = "$/unidad"
$ws.Cells.Item(893, 15).Value = "Región Metropolitana"
$ws.Cells.Item(893, 16).Value = 700
$ws.Cells.Item(893, 17).Value = 1
$ws.Cells.Item(893, 18).Value = "Hortaliza"

# Row 894: Coliflor, Segunda, Región de O'Higgins
$ws.Cells.Item(894, 1).Value = 9
$ws.Cells.Item(894, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(894, 3).Value = "Metropolitana"
$ws.Cells.Item(894, 4).Value = 45106
$ws.Cells.Item(894, 5).Value = 13
$ws.Cells.Item(894, 6).Value = 100112008
$ws.Cells.Item(894, 7).Value = "Coliflor"
$ws.Cells.Item(894, 8).Value = "Sin especificar"
$ws.Cells.Item(894, 9).Value = "Segunda"
$ws.Cells.Item(894, 10).Value = 520
$ws.Cells.Item(894, 11).Value = 800
$ws.Cells.Item(894, 12).Value = 800
$ws.Cells.Item(894, 13).Value = 800
$ws.Cells.Item(894, 14).Value = "$/unidad"
$ws.Cells.Item(894, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(894, 16).Value = 800
$ws.Cells.Item(894, 17).Value = 1
$ws.Cells.Item(894, 18).Value = "Hortaliza"
